$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (id 2) : add cohort (G) ---
$ws.Range("G3").Value = "SELECT * FROM Patient"

# --- Row 4 (id 2) : add cohort (G) ---
$ws.Range("G4").Value = "SELECT * FROM Patient"

# --- Row 5 (id 3) : add cohort (G); update data_provision (L) ---
$ws.Range("G5").Value = "SELECT * FROM Patient"
$ws.Range("L5").Value = 45573

# --- Row 6 (id 4) : add cohort (G); add data_provision (L) ---
$ws.Range("G6").Value = "SELECT * FROM Patient"
$ws.Range("L6").Value = 45478

# --- Row 7 (id 5) : add cohort (G); update data_provision (L) ---
$ws.Range("G7").Value = "SELECT * FROM Patient"
$ws.Range("L7").Value = 45534

# --- Row 8 (id 5) : add cohort (G); set date_of_submission (K) ---
$ws.Range("G8").Value = "SELECT * FROM Patient"
$ws.Range("K8").Value = 45000

# --- Row 9 (id 7) : add cohort (G); add uac_decision (J); add data_provision (L) ---
$ws.Range("G9").Value = "SELECT * FROM Patient"
$ws.Range("J9").Value = " Positive decision with conditions "
$ws.Range("L9").Value = 45414

# --- Row 10 (id 8) : add cohort (G); add uac_decision (J); update data_provision (L) ---
$ws.Range("G10").Value = "SELECT * FROM Patient"
$ws.Range("J10").Value = "Positive decision"
$ws.Range("L10").Value = 45575

# --- Row 11 (id 9) : add cohort (G); add uac_decision (J, styled like J2/J3/J4);
#     add data_provision (L); taller row to match the styled uac_decision cell ---
$ws.Range("G11").Value = "SELECT * FROM Patient"
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("J11").Value = "Categorisation of the application as not ready for decision"
$ws.Range("L11").Value = 45290
$ws.Range("A11").EntireRow.RowHeight = 15.75

# --- Row 12 (id 9) : add cohort (G); add uac_decision (J); set date_of_submission (K);
#     add data_provision (L) ---
$ws.Range("G12").Value = "SELECT * FROM Patient"
$ws.Range("J12").Value = "Positive decision"
$ws.Range("K12").Value = 45296
$ws.Range("L12").Value = 45417

# --- Row 13 (id 11) : fill in the previously-empty row entirely ---
$ws.Range("B13").Value = "Distributed analyses"
$ws.Range("C13").Value = "Manuela Musterfrau"
$ws.Range("D13").Value = "Forschungsinstitut L"
$ws.Range("E13").Value = "Bone"
$ws.Range("F13").Value = "Surgery department"
$ws.Range("G13").Value = "SELECT * FROM Patient"
$ws.Range("H13").Value = 80
$ws.Range("I13").Value = 90
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("J13").Value = "Categorisation of the application as not ready for decision"
$ws.Range("K13").Value = 45518
$ws.Range("L13").Value = 45640
$ws.Range("A13").EntireRow.RowHeight = 15.75

$excel.CutCopyMode = 0

# --- sheet selection: drop the frozen top-left scroll anchor, move the active cell ---
$ws.Range("J8").Select() | Out-Null

# --- workbook window placement ---
$win = $excel.ActiveWindow
$win.Left = -28920
$win.Top = -120
$win.Width = 29040
$win.Height = 15840
